$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Through 2022-07-09"

# Update shared string for July label (cell A8)
$ws.Range("A8").Value = "July (through 07-09)"

# Update July row (row 8) values C8:I8
$ws.Range("C8").Value = 12
$ws.Range("D8").Value = 15
$ws.Range("E8").Value = 25
$ws.Range("F8").Value = 13
$ws.Range("G8").Value = 31
$ws.Range("H8").Value = 46
$ws.Range("I8").Value = 40

# Update Total row (row 9) values C9:I9
$ws.Range("C9").Value = 260
$ws.Range("D9").Value = 405
$ws.Range("E9").Value = 378
$ws.Range("F9").Value = 264
$ws.Range("G9").Value = 503
$ws.Range("H9").Value = 806
$ws.Range("I9").Value = 846
